# Generate Report for Handoff
# Update the "5702a499-1e12-4816-b4f0-8e0f0822bf6e.md" row to reflect that the
# file is now ready for handoff (status flip) and record the new handoff
# datetime / validation error message produced by the report generator.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/51610b302e7c3fd99f3fe972a096e4c2a6906b57/e2e/5702a499-1e12-4816-b4f0-8e0f0822bf6e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bdb4ba0d8fbb5867f1805d88d7419f2fa73d6668/e2e/5702a499-1e12-4816-b4f0-8e0f0822bf6e.md."

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-04 10:53:13"

# --- zh-cn sheet ------------------------------------------------------------
# (ColumnWidth on this host stores a ~0.8333 "standard-font padding" offset
#  on top of the value you assign, so 39.1666... round-trips to the XML
#  column width of 40 that the generated report actually used.)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-09-04 10:53:08"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Range("P1:P3").ColumnWidth = 39.166666666666664

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-09-04 10:53:13"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Range("P1:P3").ColumnWidth = 39.166666666666664
